# Applies the latest cryptos-list scrape to the worksheet:
#  - refreshes Price (col D) and Volume(1h) (col E) for most rows
#  - rows 11/12 swap coin identity: TRON <-> Cardano (name, link, price, volume)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: was TRON, now Cardano
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.376"
$ws.Range("E11").Value = "  +1.84%  "

# Row 12: was Cardano, now TRON
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.156"
$ws.Range("E12").Value = "  +1.56%  "

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.334.84"
$ws.Range("E2").Value = "  -0.85%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.613.32"
$ws.Range("E3").Value = "  -0.68%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.39"
$ws.Range("E5").Value = "  +2.83%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.20"
$ws.Range("E6").Value = "  +0.15%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("E8").Value = "  -0.68%  "

# Row 9
$ws.Range("E9").Value = "  -0.13%  "

# Row 10
$ws.Range("E10").Value = "  -0.87%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.074.06"
$ws.Range("E13").Value = "  +0.09%  "

# Row 14
$ws.Range("E14").Value = "  +5.59%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.326.91"
$ws.Range("E15").Value = "  -0.87%  "

# Row 16
$ws.Range("E16").Value = "  -0.57%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.617.22"
$ws.Range("E17").Value = "  -0.13%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.45"
$ws.Range("E18").Value = "  +2.08%  "

# Row 19
$ws.Range("E19").Value = "  -0.29%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.50"
$ws.Range("E20").Value = "  -0.74%  "

# Row 21
$ws.Range("E21").Value = "  -2.43%  "

# Row 22
$ws.Range("E22").Value = "  -0.27%  "

# Row 23
$ws.Range("E23").Value = "  +2.39%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.73"
$ws.Range("E24").Value = "  -1.18%  "

# Row 25
$ws.Range("E25").Value = "  +0.25%  "

# Row 26
$ws.Range("E26").Value = "  -0.27%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.04"
$ws.Range("E27").Value = "  +4.18%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.94"
$ws.Range("E28").Value = "  +6.06%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0799"
$ws.Range("E29").Value = "  +0.37%  "

# Row 30
$ws.Range("E30").Value = "  +2.00%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.66"
$ws.Range("E31").Value = "  +4.89%  "

# Row 32
$ws.Range("E32").Value = "  +0.02%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.49"
$ws.Range("E33").Value = "  -0.25%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.03"
$ws.Range("E34").Value = "  +6.56%  "

# Row 35
$ws.Range("E35").Value = "  +0.78%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.31"
$ws.Range("E36").Value = "  +8.23%  "

# Row 37
$ws.Range("E37").Value = "  +2.79%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "318.91"
$ws.Range("E38").Value = "  +6.90%  "

# Row 40
$ws.Range("E40").Value = "  +3.55%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.851"
$ws.Range("E41").Value = "  -0.01%  "

# Row 42
$ws.Range("E42").Value = "  -3.43%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0992"
$ws.Range("E43").Value = "  +0.57%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.08%  "

# Row 45
$ws.Range("E45").Value = "  +1.76%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.608"
$ws.Range("E46").Value = "  +0.10%  "

# Row 47
$ws.Range("E47").Value = "  +3.87%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0551"
$ws.Range("E48").Value = "  +0.40%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.05"
$ws.Range("E49").Value = "  +1.54%  "

# Row 50
$ws.Range("E50").Value = "  -0.36%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.73"
$ws.Range("E51").Value = "  +0.32%  "

